$wb = $excel.ActiveWorkbook

# Delete Sheet2 entirely
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

# Work on Sheet1
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column E (engine quantizes ColumnWidth to an MDW-7 pixel grid, so use
# the input value whose quantized result lands closest to the authored 17.140625)
$ws.Columns.Item(5).ColumnWidth = 16.28

# Add dates to the "Date" column (10 Jan 2025 == serial 45667) - write the raw
# serial number so the cell keeps its existing (date) number format instead of
# Excel creating a brand-new "m/d/yyyy" style for it.
$ws.Range("E16").Value = 45667
$ws.Range("E19").Value = 45667

# Move the active selection
$ws.Range("E20").Select()
